$d = $word.ActiveDocument

# The first paragraph currently holds two runs:
#   "**ID__AFFARS_pgi_5307_topic_10__ID**" + " "
# It needs to become a single run reading
#   "**ID__AFFARS_AFMC_PGI_5307_104_92__ID**"
# and the paragraph needs a border (5pt space on all sides) and
# a deeper left indent (225 instead of 120).

$para = $d.Paragraphs(1)
$range = $para.Range

# Replace the whole paragraph's text (including the trailing space run)
# with the new placeholder text, keeping the paragraph mark intact.
$textRange = $d.Range($range.Start, $range.End - 1)
$textRange.Text = "**ID__AFFARS_AFMC_PGI_5307_104_92__ID**"

# Update paragraph formatting: indentation and border.
$para.Format.LeftIndent = 11.25
$para.Borders.DistanceFromTop = 5
$para.Borders.DistanceFromLeft = 5
$para.Borders.DistanceFromBottom = 5
$para.Borders.DistanceFromRight = 5
